$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-DataRow($row, $values) {
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($row, $i + 1).Value = $values[$i]
    }
}

# Row 2 (Toyota Yaris) keeps its values, only the name changes
$ws.Cells.Item(2, 1).Value = "Toyota Yaris  2020"

# Row 3 becomes Honda Jazz's data (previously in row 4)
Set-DataRow 3 @("Honda Jazz  2020", 17.31, 23.08, 25.64, 12.18, 5.77, 8.33, 7.69, 0, 100)

# Row 4 becomes SEAT Leon's data (previously in row 6)
Set-DataRow 4 @("SEAT Leon  2020", 0, 29.52, 24.7, 19.28, 9.039999999999999, 13.86, 3.61, 0, 100)

# Row 5 becomes Kia Sorento's data (previously in row 7)
Set-DataRow 5 @("Kia Sorento  2020", 0, 22.53, 36.26, 16.48, 9.890000000000001, 12.64, 2.2, 0, 100)

# Row 6 becomes Isuzu D-Max's data (previously in row 10)
Set-DataRow 6 @("Isuzu D-Max  2020", 0, 14.2, 51.85, 22.84, 8.640000000000001, 2.47, 0, 0, 100)

# Remove the now-unused rows 7-11 (Landrover, Honda e, Hyundai i10, old Isuzu, Audi A3)
$ws.Rows("7:11").Delete()

# Column A narrows from 21 to 20
$ws.Range("A1").ColumnWidth = 19.1
